$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.1251376819540477
$ws.Range("C2").Value = 0.2607834186411893
$ws.Range("D2").Value = -0.1914752606047823
$ws.Range("E2").Value = -0.04901477847494949
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = -0.04185119427250193
$ws.Range("H2").Value = 0.3160724455188618
$ws.Range("I2").Value = -0.3712583999854013
$ws.Range("J2").Value = 0.5861913982078005
$ws.Range("K2").Value = 0.06267724568464202
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0.03210598135678911
$ws.Range("N2").Value = 0.3028180265499283
$ws.Range("O2").Value = 0.4487784367665249
$ws.Range("P2").Value = -0.1434743588166731
$ws.Range("Q2").Value = -0.01364097849564345
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = -0.4696251170963055
$ws.Range("T2").Value = 0.2898067484975262
$ws.Range("U2").Value = -0.04848856394898793
$ws.Range("V2").Value = -0.1133556784408458

$ws.Range("B3").Value = 0.3407554847115453
$ws.Range("C3").Value = 0.04416756524540516
$ws.Range("D3").Value = 0.1427683201520291
$ws.Range("E3").Value = 0.7099644494619942
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7508659919208431
$ws.Range("H3").Value = 0.01388547601283866
$ws.Range("I3").Value = 0.003495630456440933
$ws.Range("J3").Value = 0.0000008608253554332517
$ws.Range("K3").Value = 0.6342516514928669
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8076001735331673
$ws.Range("N3").Value = 0.01868351599283946
$ws.Range("O3").Value = 0.0003226309973744454
$ws.Range("P3").Value = 0.2741132356577692
$ws.Range("Q3").Value = 0.9176101042041038
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 0.0001533968063783784
$ws.Range("T3").Value = 0.00007943445865342719
$ws.Range("U3").Value = 0.5180264530637675
$ws.Range("V3").Value = 0.1297423760635343
